$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3")
$ws.Activate()

# Fill in the newly-known frame counts for rows 15-22 (columns B and C).
# Column D already holds a shared formula (=IF(C<>"",IF(B<>"",C-B,"-"),"-"))
# that will recompute automatically once B/C are populated.
$ws.Range("B15").Value = 45855
$ws.Range("C15").Value = 55415

$ws.Range("B16").Value = 50630
$ws.Range("C16").Value = 60540

$ws.Range("B17").Value = 55302
$ws.Range("C17").Value = 65410

$ws.Range("B18").Value = 58398
$ws.Range("C18").Value = 68919

$ws.Range("B19").Value = 63243
$ws.Range("C19").Value = 74902

$ws.Range("B20").Value = 68698
$ws.Range("C20").Value = 82072

$ws.Range("B21").Value = 71616
$ws.Range("C21").Value = 85866

$ws.Range("B22").Value = 76434
$ws.Range("C22").Value = 90752

# Move the active selection to C23, matching where the user ended up after
# entering the new data.
$ws.Range("C23").Select()
